# Commit: "Added on Aug 21"
#
# - Rename the active sheet "ForgotPassword4A" -> "ForgotPassword4B"
#   (sheetId 5 / r:id rId3, stays the 3rd tab).
# - Move that sheet's selection from I25 to C15.
#
# (The absPath/revisionPtr GUID churn visible in the raw XML diff is
# Excel-internal save bookkeeping that isn't part of the scriptable
# object model - it isn't something a user/macro action sets.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ForgotPassword4B"

$ws.Activate()
$ws.Range("C15").Select()
